$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.811.36'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '2.567.04'
$ws.Range('E3').Value = '  +1.46%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.69'
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.44'
$ws.Range('E6').Value = '  +3.67%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +0.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.88'
$ws.Range('E10').Value = '  -0.70%  '
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('E12').Value = '  -0.99%  '
$ws.Range('D13').Value = '2.963.03'
$ws.Range('E13').Value = '  +1.44%  '
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.75'
$ws.Range('E15').Value = '  +3.49%  '
$ws.Range('D16').Value = '2.590.46'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.848'
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').Value = '42.868.12'
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.56'
$ws.Range('E20').Value = '  -2.15%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.38'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.07'
$ws.Range('E23').Value = '  -1.56%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.95'
$ws.Range('E27').Value = '  -0.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.38'
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.41'
$ws.Range('E29').Value = '  -0.41%  '
$ws.Range('E30').Value = '  -0.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.25'
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('E32').Value = '  -1.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.37'
$ws.Range('E33').Value = '  +0.63%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0800'
$ws.Range('E34').Value = '  +2.47%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.11'
$ws.Range('E35').Value = '  -2.81%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.76'
$ws.Range('E37').Value = '  -1.86%  '
$ws.Range('E38').Value = '  +10.16%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.118'
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.43'
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('E42').Value = '  +7.18%  '
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0302'
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.23'
$ws.Range('E45').Value = '  -2.32%  '
$ws.Range('D46').Value = '2.006.09'
$ws.Range('E46').Value = '  -1.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.96'
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('D48').Value = '2.814.34'
$ws.Range('E48').Value = '  +1.39%  '
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.66'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.65'
$ws.Range('E51').Value = '  -3.86%  '
